# considertion for LLOQ added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("F1").Value = "lloq"
$ws.Range("F2").Value = "<0.01"

# "2.23" looks numeric, so a plain assignment would store it as a number.
# Build it as a text formula result in a scratch cell, then paste-special
# values-only into F3 so it lands as a genuine shared string with no
# extra number-format/style baggage.
$scratch = $ws.Cells.Item(10, 10)
$scratch.Formula = "=""2.23"""
$scratch.Copy()
$ws.Cells.Item(3, 6).PasteSpecial(-4163)
$excel.CutCopyMode = 0
$scratch.Clear()

$ws.Activate()
$ws.Range("G11").Select()
